$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-03-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-24 Monday", 2) | Out-Null

# Update the division problems in the table, cell by cell (two cells
# share the same original text "66÷9=" but diverge to different new
# values, so positional addressing is used instead of a global replace)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "64÷8="
$tbl.Cell(1,2).Range.Text = "67÷2="
$tbl.Cell(1,3).Range.Text = "77÷6="
$tbl.Cell(1,4).Range.Text = "52÷9="
$tbl.Cell(1,5).Range.Text = "55÷9="
$tbl.Cell(5,1).Range.Text = "63÷2="
$tbl.Cell(5,2).Range.Text = "52÷8="
$tbl.Cell(5,3).Range.Text = "21÷5="
$tbl.Cell(5,4).Range.Text = "22÷2="
$tbl.Cell(5,5).Range.Text = "37÷9="
$tbl.Cell(9,1).Range.Text = "78÷8="
$tbl.Cell(9,2).Range.Text = "86÷5="
$tbl.Cell(9,3).Range.Text = "12÷2="
$tbl.Cell(9,4).Range.Text = "83÷4="
$tbl.Cell(9,5).Range.Text = "63÷9="
$tbl.Cell(13,1).Range.Text = "99÷9="
$tbl.Cell(13,2).Range.Text = "34÷8="
$tbl.Cell(13,3).Range.Text = "45÷9="
$tbl.Cell(13,4).Range.Text = "22÷2="
$tbl.Cell(13,5).Range.Text = "29÷8="
$tbl.Cell(17,1).Range.Text = "12÷5="
$tbl.Cell(17,2).Range.Text = "33÷4="
$tbl.Cell(17,3).Range.Text = "32÷3="
$tbl.Cell(17,4).Range.Text = "30÷3="
$tbl.Cell(17,5).Range.Text = "90÷6="
